$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VOCALS")

# Fill in D5 (percentage, formatted like B5/C5)
$ws.Range("D5").Value = 0.4626
$ws.Range("D5").Style = $ws.Range("C5").Style
$ws.Range("D5").NumberFormat = $ws.Range("C5").NumberFormat

# Add new "Closest" row under the EDITED row
$ws.Range("A6").Value2 = "Closest"
$ws.Range("B6").Value = 13
$ws.Range("C6").Value = 14
$ws.Range("D6").Value = 12

# Update selection to D7 as in the saved workbook
$ws.Range("D7").Select()
